$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting existing rows (and their data) down by
# one - this mirrors the author's edit: a new "dispute" entry was logged at
# 2024-09-22 20:18:07, pushing the previously-logged entries (and the
# "Broadband" label further down the sheet) down by a row.
$ws.Rows(47).Insert()

# Populate the newly inserted row with the new log entry.
$ws.Range("R47").Value = "dispute"
$ws.Range("S47").Value = "2024-09-22 20:18:07"
